# Auto-generated Excel COM-interop script
# Applies cached market-data value updates to the Excalibur_Profits workbook
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 9110.593999999999
$ws.Range("I62").Value = 7300.1665
$ws.Range("J62").Value = 11438.286
$ws.Range("K62").Value = 7300.1665
$ws.Range("L62").Value = 11438.286
$ws.Range("M62").Value = -6676.1665
$ws.Range("N62").Value = -12686.286
$ws.Range("H65").Value = 9110.593999999999
$ws.Range("I65").Value = 7300.1665
$ws.Range("J65").Value = 11438.286
$ws.Range("K65").Value = 36500.8325
$ws.Range("L65").Value = 57191.43
$ws.Range("M65").Value = -33380.8325
$ws.Range("N65").Value = -63431.43
$ws.Range("H74").Value = 8374.643
$ws.Range("I74").Value = 7926.6
$ws.Range("J74").Value = 8623.556
$ws.Range("K74").Value = 7926.6
$ws.Range("L74").Value = 8623.556
$ws.Range("M74").Value = -6990.6
$ws.Range("N74").Value = -10495.556
$ws.Range("H77").Value = 8374.643
$ws.Range("I77").Value = 7926.6
$ws.Range("J77").Value = 8623.556
$ws.Range("K77").Value = 39633
$ws.Range("L77").Value = 43117.78
$ws.Range("M77").Value = -34953
$ws.Range("N77").Value = -52477.78
$ws.Range("H98").Value = 1199.9656
$ws.Range("I98").Value = 1199.9656
$ws.Range("K98").Value = 1199.9656
$ws.Range("M98").Value = 298.0344
$ws.Range("H112").Value = 1771.8462
$ws.Range("I112").Value = 891.3333
$ws.Range("K112").Value = 2673.9999
$ws.Range("M112").Value = -1565.9999
$ws.Range("H113").Value = 1616.579
$ws.Range("J113").Value = 1281.4286
$ws.Range("L113").Value = 1281.4286
$ws.Range("N113").Value = -7789.4286
$ws.Range("H122").Value = 1199.9656
$ws.Range("I122").Value = 1199.9656
$ws.Range("K122").Value = 3599.8968
$ws.Range("M122").Value = -1149.8968
$ws.Range("H125").Value = 38463452
$ws.Range("J125").Value = 50001788
$ws.Range("L125").Value = 450016092
$ws.Range("N125").Value = -450021012
$ws.Range("H137").Value = 598144.75
$ws.Range("J137").Value = 963043.75
$ws.Range("L137").Value = 2889131.25
$ws.Range("N137").Value = -2894231.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2761.4856
$ws.Range("I74").Value = 2450.7551
$ws.Range("J74").Value = 3486.524
$ws.Range("K74").Value = 2450.7551
$ws.Range("L74").Value = 3486.524
$ws.Range("M74").Value = -1576.7551
$ws.Range("N74").Value = -5234.523999999999
$ws.Range("H77").Value = 2761.4856
$ws.Range("I77").Value = 2450.7551
$ws.Range("J77").Value = 3486.524
$ws.Range("K77").Value = 12253.7755
$ws.Range("L77").Value = 17432.62
$ws.Range("M77").Value = -7885.7755
$ws.Range("N77").Value = -26168.62
$ws.Range("H122").Value = 2638.3044
$ws.Range("I122").Value = 1392.1333
$ws.Range("K122").Value = 4176.3999
$ws.Range("M122").Value = -1726.3999
$ws.Range("H132").Value = 2967.3777
$ws.Range("I132").Value = 2475.9443
$ws.Range("J132").Value = 4933.1113
$ws.Range("K132").Value = 7427.8329
$ws.Range("L132").Value = 14799.3339
$ws.Range("M132").Value = -4897.8329
$ws.Range("N132").Value = -19859.3339

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2949.1428
$ws.Range("I20").Value = 2377.4375
$ws.Range("J20").Value = 4778.6
$ws.Range("K20").Value = 2377.4375
$ws.Range("L20").Value = 4778.6
$ws.Range("M20").Value = -2130.4375
$ws.Range("N20").Value = -5272.6

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 735.4
$ws.Range("I16").Value = 725
$ws.Range("K16").Value = 725
$ws.Range("M16").Value = -438
$ws.Range("H31").Value = 5670.482
$ws.Range("I31").Value = 2105.4827
$ws.Range("K31").Value = 2105.4827
$ws.Range("M31").Value = -1810.4827
$ws.Range("H34").Value = 5670.482
$ws.Range("I34").Value = 2105.4827
$ws.Range("K34").Value = 2105.4827
$ws.Range("M34").Value = -1903.4827
$ws.Range("H58").Value = 2497.1428
$ws.Range("I58").Value = 1838.3334
$ws.Range("K58").Value = 1838.3334
$ws.Range("M58").Value = -1635.3334
$ws.Range("H113").Value = 735.4
$ws.Range("I113").Value = 725
$ws.Range("K113").Value = 725
$ws.Range("M113").Value = 1445
$ws.Range("H122").Value = 769.5714
$ws.Range("I122").Value = 797.8333
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 2393.4999
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = 56.5001000000002
$ws.Range("N122").Value = -6700
$ws.Range("H136").Value = 2497.1428
$ws.Range("I136").Value = 1838.3334
$ws.Range("K136").Value = 5515.0002
$ws.Range("M136").Value = -2965.0002

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 9249.6
$ws.Range("I57").Value = 8549.5
$ws.Range("J57").Value = 9424.625
$ws.Range("K57").Value = 25648.5
$ws.Range("L57").Value = 28273.875
$ws.Range("M57").Value = -25089.5
$ws.Range("N57").Value = -29391.875
$ws.Range("H70").Value = 2974.75
$ws.Range("I70").Value = 2974.75
$ws.Range("K70").Value = 8924.25
$ws.Range("M70").Value = -8609.25
$ws.Range("H73").Value = 2974.75
$ws.Range("I73").Value = 2974.75
$ws.Range("K73").Value = 8924.25
$ws.Range("M73").Value = -7832.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 30529.715
$ws.Range("I13").Value = 41741.8
$ws.Range("J13").Value = 2499.5
$ws.Range("K13").Value = 41741.8
$ws.Range("L13").Value = 2499.5
$ws.Range("M13").Value = -41602.8
$ws.Range("N13").Value = -2777.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 14702.714
$ws.Range("I16").Value = 14702.714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 14702.714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -14532.714
$ws.Range("N16").ClearContents()
$ws.Range("H61").Value = 1583.0454
$ws.Range("I61").Value = 1517.3158
$ws.Range("K61").Value = 1517.3158
$ws.Range("M61").Value = -1315.3158
$ws.Range("H113").Value = 1583.0454
$ws.Range("I113").Value = 1517.3158
$ws.Range("K113").Value = 1517.3158
$ws.Range("M113").Value = 652.6841999999999
$ws.Range("H132").Value = 788709.4399999999
$ws.Range("I132").Value = 1117690.2
$ws.Range("K132").Value = 3353070.6
$ws.Range("M132").Value = -3350540.6

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2331.4285
$ws.Range("I23").Value = 466.33334
$ws.Range("K23").Value = 466.33334
$ws.Range("M23").Value = -237.33334
$ws.Range("H86").Value = 59665.668
$ws.Range("J86").Value = 59665.668
$ws.Range("L86").Value = 59665.668
$ws.Range("N86").Value = -61911.668
$ws.Range("H89").Value = 59665.668
$ws.Range("J89").Value = 59665.668
$ws.Range("L89").Value = 298328.34
$ws.Range("N89").Value = -309560.34
$ws.Range("H122").Value = 1314.8448
$ws.Range("I122").Value = 969.8182
$ws.Range("K122").Value = 2909.4546
$ws.Range("M122").Value = -459.4546
$ws.Range("H136").Value = 8470135
$ws.Range("I136").Value = 10300870
$ws.Range("J136").Value = 2983.625
$ws.Range("K136").Value = 30902610
$ws.Range("L136").Value = 8950.875
$ws.Range("M136").Value = -30900060
$ws.Range("N136").Value = -14050.875

